$p = $ppt.ActivePresentation
$s = $p.Slides.Item(17)
$shp = $s.Shapes.Item(1)

# Resize/reposition the title placeholder (shifted left, widened, same center)
$shp.Left = 335.30803149606301
$shp.Top = 242.37275590551181
$shp.Width = 289.38385826771654
$shp.Height = 55.254409448818898

# Replace the two-run text ("질문 " + "드루와") with a single run "끝"
$shp.TextFrame.TextRange.Text = "끝"
